{"js": "// Update the student name / student-ID placeholders on the title page:\n//   \u59d3\u540d\uff1aXXX   \u5b66\u53f7\uff1aXXX   ->   \u59d3\u540d\uff1a\u9ec4\u6893\u8c6a   \u5b66\u53f7\uff1a2023310103008\n//\n// The placeholders live in the second paragraph of the document\n// (\"\u59d3\u540d\uff1aXXX   \u5b66\u53f7\uff1aXXX\"). We locate that paragraph by its \"\u59d3\u540d\" label and\n// then replace the two \"XXX\" placeholders it contains, left-to-right.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet namePara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"\u59d3\u540d\") !== -1) {\n    namePara = p;\n    break;\n  }\n}\n\nif (namePara) {\n  // --- 1. \u59d3\u540d\uff1aXXX -> \u59d3\u540d\uff1a\u9ec4\u6893\u8c6a ------------------------------------------\n  let placeholders = namePara.search(\"XXX\", { matchCase: true });\n  placeholders.load(\"items\");\n  await context.sync();\n\n  if (placeholders.items.length > 0) {\n    placeholders.items[0].insertText(\"\u9ec4\u6893\u8c6a\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // --- 2. \u5b66\u53f7\uff1aXXX -> \u5b66\u53f7\uff1a2023310103008 -----------------------------------\n  // Only the student-id placeholder is left now.\n  placeholders = namePara.search(\"XXX\", { matchCase: true });\n  placeholders.load(\"items\");\n  await context.sync();\n\n  if (placeholders.items.length > 0) {\n    placeholders.items[0].insertText(\"2023310103008\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // Shift the \"\uff1a\" from the end of the \"\u5b66\u53f7\uff1a\" label onto the id run so the\n  // run boundaries match the final document (\u5b66\u53f7 | \uff1a2023310103008).\n  let label = namePara.search(\"\u5b66\u53f7\uff1a\", { matchCase: true });\n  label.load(\"items\");\n  await context.sync();\n\n  if (label.items.length > 0) {\n    label.items[0].insertText(\"\u5b66\u53f7\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  let idNumber = namePara.search(\"2023310103008\", { matchCase: true });\n  idNumber.load(\"items\");\n  await context.sync();\n\n  if (idNumber.items.length > 0) {\n    idNumber.items[0].insertText(\"\uff1a2023310103008\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue=1, wdReplaceOne=1 -> replace only the first match found.\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n\n# Title page placeholders:\n#   \u59d3\u540d\uff1aXXX   \u5b66\u53f7\uff1aXXX   ->   \u59d3\u540d\uff1a\u9ec4\u6893\u8c6a   \u5b66\u53f7\uff1a2023310103008\n# Replace the name placeholder first (leftmost \"XXX\"), then the id\n# placeholder (the \"XXX\" that remains after the first replace).\nReplace-FirstMatch \"XXX\" \"\u9ec4\u6893\u8c6a\"\nReplace-FirstMatch \"XXX\" \"2023310103008\"\n"}
